$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column C width (new custom width added in diff)
$ws.Columns.Item(3).ColumnWidth = 21

# Update balance value for row 5 (E5: 2700 -> 600)
$ws.Range("E5").Value = 600

# Move the active selection from B10 to E9
$ws.Range("E9").Select()
